$wb = $excel.ActiveWorkbook

# "展览" sheet (1st sheet / sheetId 1)
$wsExhibit = $wb.Worksheets.Item(1)
$wsExhibit.Range("F2").Value = 1080
$wsExhibit.Range("F4").Value = 1646
$wsExhibit.Range("F5").Value = 749
$wsExhibit.Range("F6").Value = 88

# "全部类型" sheet (4th sheet / sheetId 4)
$wsAll = $wb.Worksheets.Item(4)
$wsAll.Range("F2").Value = 1080
$wsAll.Range("F4").Value = 1646
$wsAll.Range("F6").Value = 749
$wsAll.Range("F7").Value = 88
